$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Session 7 (Branch&Bound) column for student "Cuesta Martínez, Miguel": mark 0
$ws.Range("H4").Value = 0

# Teacher's feedback for session 7
$ws.Range("H5").Value = "Not implemented"

# Update the active selection to match new focus cell H5:H12
$ws.Range("H5:H12").Select()
